$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the merge-diagnostic columns to their new "_mp_" prefixed names
# (the library that produced this file switched from a bare "_merge" /
# "_diff_days" / "_abs_diff_days" naming scheme to a "_mp_"-prefixed one).
$ws.Range("O1").Value = "_mp_merge"
$ws.Range("P1").Value = "_mp_diff_days"
$ws.Range("Q1").Value = "_mp_abs_diff_days"

# Only rows that actually matched in the merge ("both") keep a value in
# the _mp_merge column now; the "left_only" marker cells are cleared.
$ws.Range("O2").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("O7").ClearContents()
$ws.Range("O8").ClearContents()
$ws.Range("O9").ClearContents()
$ws.Range("O11").ClearContents()
$ws.Range("O12").ClearContents()
$ws.Range("O13").ClearContents()

# The whole "_duplicates" column (R) no longer exists in the refactored
# output, so remove it entirely (this shifts nothing else, R was last).
$ws.Columns.Item(18).Delete()

# Approximate the best-fit column widths Excel computed for the
# newly-widened, renamed headers.
$ws.Columns.Item(15).ColumnWidth = 9.498697916666666
$ws.Columns.Item(16).ColumnWidth = 11.166666666666666
$ws.Columns.Item(17).ColumnWidth = 14.498697916666666
